# The commit re-shuffles the content of rows 44-53 on the "Artfynd" sheet:
# each whole row's data moves to a different row within the same 44-53
# block (a permutation), while row/column structure, headers and all other
# rows stay untouched.
#
# Mapping: destination row -> source row (i.e. row 44 ends up holding what
# used to be row 53's data, row 45 ends up holding what used to be row
# 51's data, etc). Derived from the authoritative diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowMap = @{
    44 = 53
    45 = 51
    46 = 48
    47 = 46
    48 = 47
    49 = 52
    50 = 45
    51 = 44
    52 = 50
    53 = 49
}

$firstCol = 1    # A
$lastCol  = 51   # AY

# Columns whose text content could otherwise be auto-coerced by Excel on
# plain assignment (pure numeric strings in "Antal", ISO dates in
# "Startdatum"/"Slutdatum") - these must stay text, so we force the "@"
# (Text) number format on them before writing.
$textForceCols = @(9, 25, 27)   # I=Antal, Y=Startdatum, AA=Slutdatum

# --- Step 1: snapshot every source cell (rows 44-53, cols A:AY) before
# any writes happen, so overlapping source/destination rows in the
# permutation don't clobber data we still need to read. ---
$snapshot = @{}
foreach ($srcRow in $rowMap.Values) {
    if (-not $snapshot.ContainsKey($srcRow)) {
        $rowVals = @{}
        for ($c = $firstCol; $c -le $lastCol; $c++) {
            $rowVals[$c] = $ws.Cells.Item($srcRow, $c).Value2
        }
        $snapshot[$srcRow] = $rowVals
    }
}

# --- Step 2: write the snapshotted data into its new destination row. ---
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $rowVals = $snapshot[$srcRow]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($destRow, $c)
        if ($textForceCols -contains $c) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $rowVals[$c]
    }
}

Write-Host "Reshuffled rows 44-53 per the commit's row permutation."
